# Project 2 - code complete, working on report
#
# Update the "Testing" sheet's Correct/Incorrect tallies for "Basic
# Problems C" (row 2): Correct 7 -> 8, Incorrect 4 -> 3. The dependent
# SUM()/ratio formulas in rows 4-5 recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testing")
$ws.Activate()

$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 3

$ws.Range("C3").Select() | Out-Null
